# "Clean System Level Clean System Wide"
#
# 1. Move the "emission_system" tab so it sits after "emission" (i.e. right
#    before "material_cost") instead of right after "baseline".
# 2. Replace the "global" (system-wide) row on the "emission" sheet with a
#    flat cap of 999,999,999 for every year.
# 3. Replace the per-plant ("system level") rows on "emission_system" with
#    a new ramp-down schedule.
# 4. Add a "min" column (alongside the existing "max" column) to the
#    technology_fuel_pairs and technology_material_pairs lookup tables, and
#    refresh/extend their data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Re-order the sheet tabs: move emission_system after emission
# ---------------------------------------------------------------------
$wsSystem = $wb.Worksheets.Item("emission_system")
$wsAfter  = $wb.Worksheets.Item("emission")
$wsSystem.Move($null, $wsAfter)

# ---------------------------------------------------------------------
# 2) "emission" sheet - flatten the "global" row to 999999999 everywhere
# ---------------------------------------------------------------------
$wsEmission = $wb.Worksheets.Item("emission")
for ($c = 2; $c -le 27; $c++) {
    $wsEmission.Cells.Item(2, $c).Value = 999999999
}

# ---------------------------------------------------------------------
# 3) "emission_system" sheet - new per-plant schedules
# ---------------------------------------------------------------------
$wsSys = $wb.Worksheets.Item("emission_system")

# Rows 2 & 3 (POSCO1 / POSCO2): flat 1e7 through col G (2030), halve twice
# via formulas, then flat at the resulting 2.5e6 through the end.
foreach ($r in 2, 3) {
    $wsSys.Cells.Item($r, 2).Value = 99999999
    for ($c = 3; $c -le 7; $c++) {
        $wsSys.Cells.Item($r, $c).Value = 10000000
    }
    $wsSys.Cells.Item($r, 8).Formula = "=" + $wsSys.Cells.Item($r, 7).Address($false, $false) + "/2"
    $wsSys.Cells.Item($r, 9).Formula = "=" + $wsSys.Cells.Item($r, 8).Address($false, $false) + "/2"
    for ($c = 10; $c -le 27; $c++) {
        $wsSys.Cells.Item($r, $c).Formula = "=" + $wsSys.Cells.Item($r, $c - 1).Address($false, $false)
    }
}

# Rows 4 & 5 (HYUNDAI1 / HYUNDAI2): flat 1e7 through col G (2030), then
# drop to a flat 500000 for the remainder.
foreach ($r in 4, 5) {
    $wsSys.Cells.Item($r, 2).Value = 99999999
    for ($c = 3; $c -le 7; $c++) {
        $wsSys.Cells.Item($r, $c).Value = 10000000
    }
    for ($c = 8; $c -le 27; $c++) {
        $wsSys.Cells.Item($r, $c).Value = 500000
    }
}

# ---------------------------------------------------------------------
# 4) technology_fuel_pairs - add "min" column, refresh max/min values
# ---------------------------------------------------------------------
$wsFuel = $wb.Worksheets.Item("technology_fuel_pairs")
$wsFuel.Range("D1").Value = "min"

$fuelMax = @(1, 0.5, 0.4, 1, 0.2, 1, 0.2, 1, 0.5)
$fuelMin = @(0.8, 0.2, 0, 0.5, 0, 0.5, 0, 0.8, 0)
for ($i = 0; $i -lt 9; $i++) {
    $r = $i + 2
    $wsFuel.Cells.Item($r, 3).Value = $fuelMax[$i]
    $wsFuel.Cells.Item($r, 4).Value = $fuelMin[$i]
}

# ---------------------------------------------------------------------
# 5) technology_material_pairs - add "min" column, new pairing table
# ---------------------------------------------------------------------
$wsMat = $wb.Worksheets.Item("technology_material_pairs")
$wsMat.Range("D1").Value = "min"

$matTech  = @("BF-BOF", "BF-BOF", "DRI-EAF", "DRI-EAF", "DRI-EAF", "ESF", "BF-BOF-CCUS", "BF-BOF-CCUS")
$matMat   = @("Cokes", "Scrap", "NG-DRI", "H2-DRI", "Scrap", "Scrap", "Cokes", "Scrap")
$matMax   = @(1, 0.5, 1, 1, 0.2, 1, 1, 0.2)
$matMin   = @(0.2, 0.2, 0, 0, 0, 0, 0.5, 0)
for ($i = 0; $i -lt 8; $i++) {
    $r = $i + 2
    $wsMat.Cells.Item($r, 1).Value = $matTech[$i]
    $wsMat.Cells.Item($r, 2).Value = $matMat[$i]
    $wsMat.Cells.Item($r, 3).Value = $matMax[$i]
    $wsMat.Cells.Item($r, 4).Value = $matMin[$i]
}
